$d = $word.ActiveDocument

$d.Content.Find.Execute("350,000", $true, $false, $false, $false, $false, $true, 1, $false, "500,000", 2)
$d.Content.Find.Execute("defence", $true, $false, $false, $false, $false, $true, 1, $false, "defense", 2)
$d.Content.Find.Execute("assist the Soviet Logistics", $true, $false, $false, $false, $false, $true, 1, $false, "assist Soviet Logistics", 2)
